$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking strings in columns D and E stay as text, matching
# the original inlineStr cell type (avoids Excel auto-converting to numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.661.51"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.97%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.774.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.80%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "355.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "109.14"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -2.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.552"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.44%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -3.38%  "

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.00%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0842"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.40%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.44"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.70%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.73%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.210.24"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.759.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.928"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.55%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.595.40"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.66%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.44"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.50%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.34%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.24%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.98"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.30%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "267.40"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.83%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.56%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.04%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +13.38%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.85%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.71%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "51.43"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.89%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "34.36"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.31%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0444"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -9.26%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0835"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.16"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -8.31%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.88"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.75%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -5.40%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.96"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -4.54%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.68%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.22"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -2.62%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "119.37"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -5.94%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.76"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -6.46%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.088.19"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.11%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.26"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.66%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.935"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.24%  "

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.98%  "

$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.13%  "
